# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" (used by the Slide
#                            Master, i.e. the theme that is active for
#                            every slide in Normal view)
#
# The authored change swaps the two themes' content, so the deck's
# visible/active theme (theme2.xml, reached from the Slide Master) becomes
# the plain "Office Theme" colour scheme that used to live only in
# theme1.xml. Re-create that swap by rewriting the twelve theme colours on
# the Slide Master's theme to the "Office Theme" values.

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target colours, in MsoThemeColorSchemeIndex order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) = "Office Theme".
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgb($officeThemeColors[$i - 1])
}
